$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.72194870874295
$ws.Range("C2").Value = 0.163891712200428
$ws.Range("B3").Value = 0.255836150826446
$ws.Range("C3").Value = 0.118697956862633

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.29516138426015
$ws.Range("C2").Value = 0.222364136608754
$ws.Range("B3").Value = -1.17781611095422
$ws.Range("C3").Value = 0.104923760248021

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.97028949193194
$ws.Range("C2").Value = 0.0801902518794047
$ws.Range("B3").Value = 2.12084194642728
$ws.Range("C3").Value = 0.213035641702335

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.29936955689108
$ws.Range("C2").Value = 0.105835268383319
$ws.Range("B3").Value = 0.00178324955458855
$ws.Range("C3").Value = 0.0170697952170677

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.026860493327988
$ws.Range("B2").Value = -0.0149839404293542
$ws.Range("A3").Value = -0.0149839404293542
$ws.Range("B3").Value = 0.0140892049633635

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0494458092497564
$ws.Range("B2").Value = -0.0206429316639446
$ws.Range("A3").Value = -0.0206429316639446
$ws.Range("B3").Value = 0.0110089954645841

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00643047649648237
$ws.Range("B2").Value = -0.00257506918792403
$ws.Range("A3").Value = -0.00257506918792403
$ws.Range("B3").Value = 0.0453841846355257

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0112011040337692
$ws.Range("B2").Value = -0.0010767979983474
$ws.Range("A3").Value = -0.0010767979983474
$ws.Range("B3").Value = 0.000291377908752626
